$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45:A87").FormulaArray = "=LOWER(A2:A44)"
$ws.Range("B45:B87").FormulaArray = "=LOWER(B2:B44)"

[void]$ws.Range("B1").Select()
